$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7279
$ws.Range("J3").Value = 7666
$ws.Range("I4").Value = 1715
$ws.Range("J4").Value = 1667
$ws.Range("J5").Value = 599
$ws.Range("J6").Value = 10446
$ws.Range("I7").Value = 24973
$ws.Range("J7").Value = 27657

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J6").Value = 212
$ws.Range("J7").Value = 791
$ws.Range("J8").Value = 1742
$ws.Range("J10").Value = 202
$ws.Range("J11").Value = 497
$ws.Range("J13").Value = 34
$ws.Range("J17").Value = 37
$ws.Range("J19").Value = 794
$ws.Range("J20").Value = 591
$ws.Range("J23").Value = 254
$ws.Range("J29").Value = 1474
$ws.Range("J33").Value = 1254
$ws.Range("J37").Value = 847
$ws.Range("J42").Value = 1185
$ws.Range("J50").Value = 163
$ws.Range("J51").Value = 341
$ws.Range("J52").Value = 708
$ws.Range("J53").Value = 415
$ws.Range("J57").Value = 131
$ws.Range("J59").Value = 32
$ws.Range("J60").Value = 163
$ws.Range("I63").Value = 181
$ws.Range("J63").Value = 78
$ws.Range("J64").Value = 182
$ws.Range("J65").Value = 695
$ws.Range("J66").Value = 83
$ws.Range("J67").Value = 1023
$ws.Range("J71").Value = 90
$ws.Range("J75").Value = 82
$ws.Range("J76").Value = 395
$ws.Range("J77").Value = 192
$ws.Range("J79").Value = 761
$ws.Range("J83").Value = 553
$ws.Range("J84").Value = 229
$ws.Range("J85").Value = 1134
$ws.Range("J87").Value = 91
$ws.Range("J90").Value = 291
$ws.Range("J91").Value = 316
$ws.Range("J94").Value = 304
$ws.Range("J95").Value = 396
$ws.Range("J97").Value = 253
$ws.Range("J99").Value = 424
$ws.Range("I101").Value = 24973
$ws.Range("J101").Value = 27657

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 254
$ws.Range("J7").Value = 791

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 236
$ws.Range("J7").Value = 497

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 304
$ws.Range("J3").Value = 408
$ws.Range("J6").Value = 323
$ws.Range("J7").Value = 1134

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 166
$ws.Range("J6").Value = 306
$ws.Range("J7").Value = 708

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 76
$ws.Range("J7").Value = 415

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 459
$ws.Range("J3").Value = 507
$ws.Range("J6").Value = 641
$ws.Range("J7").Value = 1742

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J6").Value = 148
$ws.Range("J7").Value = 553

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 283
$ws.Range("J3").Value = 416
$ws.Range("J6").Value = 447
$ws.Range("J7").Value = 1254

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 139
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 396

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 254
$ws.Range("J3").Value = 285
$ws.Range("J6").Value = 246
$ws.Range("J7").Value = 847

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J5").Value = 19
$ws.Range("J6").Value = 259
$ws.Range("J7").Value = 695

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 113
$ws.Range("J7").Value = 424

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 262
$ws.Range("J7").Value = 1023

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 74
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 522
$ws.Range("J4").Value = 80
$ws.Range("J7").Value = 1474

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value = 306
$ws.Range("J7").Value = 794

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 71
$ws.Range("J3").Value = 87
$ws.Range("J7").Value = 395

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 246
$ws.Range("J6").Value = 627
$ws.Range("J7").Value = 1185

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("J3").Value = 11
$ws.Range("J6").Value = 34

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J3").Value = 34
$ws.Range("J7").Value = 202

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 130
$ws.Range("J7").Value = 316

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 215
$ws.Range("J3").Value = 254
$ws.Range("J6").Value = 229
$ws.Range("J7").Value = 761

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 45
$ws.Range("J7").Value = 591

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 160
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J2").Value = 42
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 175
$ws.Range("J7").Value = 253

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 67
$ws.Range("J6").Value = 153

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 80
$ws.Range("J7").Value = 291

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J4").Value = 31
$ws.Range("J7").Value = 341

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 56
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 90

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 192

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 91
